$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while preserving its original (default) style,
# forcing text interpretation for strings that look like numbers so that
# e.g. "239.09" or "2.00" are not silently coerced into numeric values.
function Set-TextValue($ws, $cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "35.144.41"
$ws.Range("E2").Value = "  +1.21%  "
$ws.Range("D3").Value = "1.856.42"
$ws.Range("E3").Value = "  +2.02%  "
$ws.Range("E4").Value = "  +0.17%  "
Set-TextValue $ws "D5" "239.09"
$ws.Range("E5").Value = "  +3.78%  "
$ws.Range("E6").Value = "  +1.43%  "
Set-TextValue $ws "D8" "41.82"
$ws.Range("E8").Value = "  +5.86%  "
$ws.Range("E9").Value = "  +3.04%  "
$ws.Range("E10").Value = "  +1.67%  "
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("D12").Value = "2.126.43"
$ws.Range("E12").Value = "  +1.98%  "
Set-TextValue $ws "D13" "11.46"
$ws.Range("E13").Value = "  +1.89%  "
$ws.Range("D14").Value = "1.849.49"
$ws.Range("E14").Value = "  +1.76%  "
Set-TextValue $ws "D15" "0.676"
$ws.Range("E15").Value = "  +1.86%  "
Set-TextValue $ws "D16" "4.71"
$ws.Range("E16").Value = "  +2.03%  "
$ws.Range("D17").Value = "35.133.04"
$ws.Range("E17").Value = "  +1.39%  "
Set-TextValue $ws "D18" "69.71"
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("E19").Value = "  +1.60%  "
Set-TextValue $ws "D20" "240.56"
Set-TextValue $ws "D21" "12.22"
$ws.Range("E21").Value = "  +1.79%  "
$ws.Range("E22").Value = "  +1.99%  "
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("E24").Value = "  +0.25%  "
Set-TextValue $ws "D25" "168.87"
$ws.Range("E25").Value = "  -2.36%  "
$ws.Range("E26").Value = "  +28.01%  "
Set-TextValue $ws "D27" "7.97"
$ws.Range("E27").Value = "  +3.73%  "
$ws.Range("E28").Value = "  +2.29%  "
$ws.Range("E29").Value = "  +0.40%  "
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("E31").Value = "  +1.88%  "
$ws.Range("E32").Value = "  +2.51%  "
Set-TextValue $ws "D33" "1.83"
$ws.Range("E33").Value = "  +27.74%  "
$ws.Range("E34").Value = "  +2.72%  "
Set-TextValue $ws "D35" "0.826"
$ws.Range("E35").Value = "  +18.57%  "
Set-TextValue $ws "D36" "2.00"
$ws.Range("E36").Value = "  +10.50%  "
$ws.Range("E37").Value = "  +7.16%  "
$ws.Range("E38").Value = "  +7.76%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws "D39" "0.0201"
$ws.Range("E39").Value = "  +4.45%  "
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws "D40" "89.64"
$ws.Range("E40").Value = "  -1.62%  "
$ws.Range("D41").Value = "1.342.07"
$ws.Range("E41").Value = "  +0.32%  "
Set-TextValue $ws "D42" "14.93"
$ws.Range("E42").Value = "  +3.86%  "
$ws.Range("E43").Value = "  +4.06%  "
$ws.Range("E44").Value = "  -0.61%  "
Set-TextValue $ws "D45" "12.38"
$ws.Range("E45").Value = "  +48.07%  "
$ws.Range("E46").Value = "  +6.83%  "
$ws.Range("E47").Value = "  +0.11%  "
Set-TextValue $ws "D48" "6.60"
$ws.Range("E48").Value = "  +5.73%  "
$ws.Range("D49").Value = "2.038.24"
$ws.Range("E49").Value = "  +1.84%  "
$ws.Range("E50").Value = "  +1.91%  "
$ws.Range("E51").Value = "  +0.16%  "
